$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.413.34"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "3.505.79"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'591.31"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").Value = "'134.73"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("E9").Value = "  +6.22%  "

$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").Value = "'0.392"
$ws.Range("E11").Value = "  +4.46%  "

$ws.Range("D12").Value = "4.102.37"
$ws.Range("E12").Value = "  +0.55%  "

$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("E14").Value = "  +1.41%  "

$ws.Range("D15").Value = "3.506.01"

$ws.Range("D16").Value = "'25.84"
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("D17").Value = "64.408.25"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'10.07"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("D19").Value = "'5.78"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "'392.02"
$ws.Range("E21").Value = "  +2.17%  "

$ws.Range("D22").Value = "'0.585"
$ws.Range("E22").Value = "  +3.53%  "

$ws.Range("D23").Value = "3.645.63"
$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("D24").Value = "'74.50"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("E27").Value = "  +4.74%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.47"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("E30").Value = "  +2.25%  "

$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").Value = "'0.158"
$ws.Range("E33").Value = "  +7.59%  "

$ws.Range("D34").Value = "3.533.82"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'23.47"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("D37").Value = "'5.37"
$ws.Range("E37").Value = "  +2.14%  "

$ws.Range("D38").Value = "'6.99"
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("D39").Value = "'1.58"
$ws.Range("E39").Value = "  +3.18%  "

$ws.Range("D40").Value = "'165.63"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("D41").Value = "'0.0794"
$ws.Range("E41").Value = "  +2.06%  "

$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").Value = "'4.47"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").Value = "'24.96"
$ws.Range("E45").Value = "  -1.84%  "

$ws.Range("D46").Value = "'1.19"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").Value = "'1.67"
$ws.Range("E47").Value = "  +1.97%  "

$ws.Range("D48").Value = "'0.928"
$ws.Range("E48").Value = "  +3.85%  "

$ws.Range("D49").Value = "'6.84"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").Value = "2.405.74"
$ws.Range("E50").Value = "  -2.43%  "

$ws.Range("E51").Value = "  +0.93%  "
